$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Shift R/S (subject/date) column values down by one row for rows 40-144,
# inserting a new top entry at row 40 ("balance your axis" / 2024-09-16 07:57:00).
$ws.Range("R40").Value = 'balance your axis'
$ws.Range("S40").Value = '2024-09-16 07:57:00'
$ws.Range("R41").Value = 'money google icici'
$ws.Range("S41").Value = '2024-09-15 21:06:00'
$ws.Range("R42").Value = 'adani icici'
$ws.Range("S42").Value = '2024-09-15 13:10:50'
$ws.Range("R43").Value = 'balance your axis'
$ws.Range("S43").Value = '2024-09-15 07:56:24'
$ws.Range("R44").Value = 'bal axisbank w axis'
$ws.Range("S44").Value = '2024-09-15 07:12:01'
$ws.Range("R45").Value = 'hdfc'
$ws.Range("S45").Value = '2024-09-14 21:25:23'
$ws.Range("R46").Value = 'change the'
$ws.Range("S46").Value = '2024-09-12 21:16:38'
$ws.Range("R47").Value = 'dispute'
$ws.Range("S47").Value = '2024-09-12 19:02:14'
$ws.Range("R48").Value = 'congrats limit icici'
$ws.Range("S48").Value = '2024-09-12 19:03:39'
$ws.Range("R49").Value = 'latest transaction pan'
$ws.Range("S49").Value = '2024-09-12 12:22:12'
$ws.Range("R50").Value = 'assistance'
$ws.Range("S50").Value = '2024-09-12 12:17:33'
$ws.Range("R51").Value = 'balance your axis'
$ws.Range("S51").Value = '2024-09-12 09:37:28'
$ws.Range("R52").Value = 'bal axisbank'
$ws.Range("S52").Value = '2024-09-12 00:54:39'
$ws.Range("R53").Value = 'your relationship'
$ws.Range("S53").Value = '2024-09-11 16:05:27'
$ws.Range("R54").Value = 'login internet personal share'
$ws.Range("S54").Value = '2024-09-11 14:16:45'
$ws.Range("R55").Value = 'balance your axis'
$ws.Range("S55").Value = '2024-09-11 12:45:33'
$ws.Range("R56").Value = 'balance your axis'
$ws.Range("S56").Value = '2024-09-11 09:45:01'
$ws.Range("R57").Value = 'axis'
$ws.Range("S57").Value = '2024-09-11 06:57:42'
$ws.Range("R58").Value = 'money google icici'
$ws.Range("S58").Value = '2024-09-10 20:42:12'
$ws.Range("R59").Value = 'dispute'
$ws.Range("S59").Value = '2024-09-10 20:42:34'
$ws.Range("R60").Value = 'reward points cash'
$ws.Range("S60").Value = '2024-09-10 19:43:35'
$ws.Range("R61").Value = 'balance your axis'
$ws.Range("S61").Value = '2024-09-10 13:32:42'
$ws.Range("R62").Value = 'ach indianesign bal axisbank'
$ws.Range("S62").Value = '2024-09-10 13:22:37'
$ws.Range("R63").Value = 'ach indianesign bal axisbank'
$ws.Range("S63").Value = '2024-09-10 13:22:37'
$ws.Range("R64").Value = 'balance your axis'
$ws.Range("S64").Value = '2024-09-10 11:21:40'
$ws.Range("R65").Value = 'your relationship'
$ws.Range("S65").Value = '2024-09-10 11:02:23'
$ws.Range("R66").Value = 'bank bal broker'
$ws.Range("S66").Value = '2024-09-09 19:59:02'
$ws.Range("R67").Value = 'beneficiary'
$ws.Range("S67").Value = '2024-09-09 15:48:10'
$ws.Range("R68").Value = 'beneficiary saravanan'
$ws.Range("S68").Value = '2024-09-09 14:52:20'
$ws.Range("R69").Value = 'bal axisbank'
$ws.Range("S69").Value = '2024-09-09 12:19:34'
$ws.Range("R70").Value = 'bal axisbank'
$ws.Range("S70").Value = '2024-09-09 12:19:33'
$ws.Range("R71").Value = 'dispute'
$ws.Range("S71").Value = '2024-09-09 12:17:30'
$ws.Range("R72").Value = 'bal axisbank'
$ws.Range("S72").Value = '2024-09-09 12:04:31'
$ws.Range("R73").Value = 'transfer freedom share anyone axis'
$ws.Range("S73").Value = '2024-09-09 11:56:19'
$ws.Range("R74").Value = 'corporate internet share'
$ws.Range("S74").Value = '2024-09-09 11:40:49'
$ws.Range("R75").Value = 'corporate internet share'
$ws.Range("S75").Value = '2024-09-09 11:39:30'
$ws.Range("R76").Value = 'bal axisbank'
$ws.Range("S76").Value = '2024-09-09 11:38:16'
$ws.Range("R77").Value = 'bal axisbank'
$ws.Range("S77").Value = '2024-09-09 11:38:16'
$ws.Range("R78").Value = 'bal axisbank'
$ws.Range("S78").Value = '2024-09-09 11:38:15'
$ws.Range("R79").Value = 'bal axisbank'
$ws.Range("S79").Value = '2024-09-09 11:38:15'
$ws.Range("R80").Value = 'corporate internet share'
$ws.Range("S80").Value = '2024-09-09 11:35:34'
$ws.Range("R81").Value = 'corporate internet share'
$ws.Range("S81").Value = '2024-09-09 11:32:23'
$ws.Range("R82").Value = 'ift anbu tpar'
$ws.Range("S82").Value = '2024-09-09 11:27:52'
$ws.Range("R83").Value = 'balance your axis'
$ws.Range("S83").Value = '2024-09-09 11:24:00'
$ws.Range("R84").Value = 'corporate internet share'
$ws.Range("S84").Value = '2024-09-09 11:21:43'
$ws.Range("R85").Value = 'corporate internet share'
$ws.Range("S85").Value = '2024-09-09 11:17:34'
$ws.Range("R86").Value = 'corporate internet share'
$ws.Range("S86").Value = '2024-09-09 11:15:51'
$ws.Range("R87").Value = 'corporate internet share'
$ws.Range("S87").Value = '2024-09-09 11:14:13'
$ws.Range("R88").Value = 'anbu tparty bal axisbank'
$ws.Range("S88").Value = '2024-09-09 11:13:37'
$ws.Range("R89").Value = 'corporate internet share'
$ws.Range("S89").Value = '2024-09-09 11:10:39'
$ws.Range("R90").Value = 'corporate internet share'
$ws.Range("S90").Value = '2024-09-09 11:07:31'
$ws.Range("R91").Value = 'corporate internet share'
$ws.Range("S91").Value = '2024-09-09 11:03:09'
$ws.Range("R92").Value = 'saravanan'
$ws.Range("S92").Value = '2024-09-09 10:43:11'
$ws.Range("R93").Value = 'balance your axis'
$ws.Range("S93").Value = '2024-09-09 08:10:16'
$ws.Range("R94").Value = 'ekalaivan'
$ws.Range("S94").Value = '2024-09-08 18:40:34'
$ws.Range("R95").Value = 'balance your axis'
$ws.Range("S95").Value = '2024-09-08 09:53:37'
$ws.Range("R96").Value = 'balance your axis'
$ws.Range("S96").Value = '2024-09-07 12:12:22'
$ws.Range("R97").Value = 'balance your axis'
$ws.Range("S97").Value = '2024-09-07 09:34:58'
$ws.Range("R98").Value = 'bal axis'
$ws.Range("S98").Value = '2024-09-07 08:46:40'
$ws.Range("R99").Value = 'axis'
$ws.Range("S99").Value = '2024-09-07 08:31:28'
$ws.Range("R100").Value = 'your relationship'
$ws.Range("S100").Value = '2024-09-06 12:23:25'
$ws.Range("R101").Value = 'balance your axis'
$ws.Range("S101").Value = '2024-09-06 09:55:31'
$ws.Range("R102").Value = 'beneficiary'
$ws.Range("S102").Value = '2024-09-05 17:13:56'
$ws.Range("R103").Value = 'coimbatore ramalinga'
$ws.Range("S103").Value = '2024-09-05 17:06:01'
$ws.Range("R104").Value = 'beneficiary'
$ws.Range("S104").Value = '2024-09-05 17:04:10'
$ws.Range("R105").Value = 'bal axisbank'
$ws.Range("S105").Value = '2024-09-05 16:52:25'
$ws.Range("R106").Value = 'share anyone axis'
$ws.Range("S106").Value = '2024-09-05 16:38:59'
$ws.Range("R107").Value = 'transfer anyone axis'
$ws.Range("S107").Value = '2024-09-05 16:35:58'
$ws.Range("R108").Value = 'share anyone axis'
$ws.Range("S108").Value = '2024-09-05 16:31:34'
$ws.Range("R109").Value = 'transfer'
$ws.Range("S109").Value = '2024-09-05 16:28:38'
$ws.Range("R110").Value = 'bal axisbank axis'
$ws.Range("S110").Value = '2024-09-05 16:26:56'
$ws.Range("R111").Value = 'bal axisbank'
$ws.Range("S111").Value = '2024-09-05 16:26:55'
$ws.Range("R112").Value = 'transfer'
$ws.Range("S112").Value = '2024-09-05 16:25:07'
$ws.Range("R113").Value = 'transfer'
$ws.Range("S113").Value = '2024-09-05 16:22:23'
$ws.Range("R114").Value = 'share anyone axis'
$ws.Range("S114").Value = '2024-09-05 16:06:05'
$ws.Range("R115").Value = 'internet bal axisbank'
$ws.Range("S115").Value = '2024-09-05 16:05:55'
$ws.Range("R116").Value = 'transfer share anyone axis'
$ws.Range("S116").Value = '2024-09-05 16:03:14'
$ws.Range("R117").Value = 'axis'
$ws.Range("S117").Value = '2024-09-05 15:57:15'
$ws.Range("R118").Value = 'your net internet'
$ws.Range("S118").Value = '2024-09-05 15:57:15'
$ws.Range("R119").Value = 'hear your feedback atm'
$ws.Range("S119").Value = '2024-09-05 14:21:08'
$ws.Range("R120").Value = 'axis bna'
$ws.Range("S120").Value = '2024-09-05 14:18:32'
$ws.Range("R121").Value = 'axis bna'
$ws.Range("S121").Value = '2024-09-05 14:13:16'
$ws.Range("R122").Value = 'axis bna'
$ws.Range("S122").Value = '2024-09-05 14:15:23'
$ws.Range("R123").Value = 'balance your axis'
$ws.Range("S123").Value = '2024-09-05 09:20:57'
$ws.Range("R124").Value = 'bal axis'
$ws.Range("S124").Value = '2024-09-05 09:06:25'
$ws.Range("R125").Value = 'broker'
$ws.Range("S125").Value = '2024-09-04 21:20:47'
$ws.Range("R126").Value = 'exclusive on axis'
$ws.Range("S126").Value = '2024-09-04 13:21:05'
$ws.Range("R127").Value = 'your corporate axis'
$ws.Range("S127").Value = '2024-09-04 11:46:10'
$ws.Range("R128").Value = 'balance your axis'
$ws.Range("S128").Value = '2024-09-04 08:14:16'
$ws.Range("R129").Value = 'axis'
$ws.Range("S129").Value = '2024-09-04 07:02:13'
$ws.Range("R130").Value = 'bal axisbank w axis'
$ws.Range("S130").Value = '2024-09-04 06:53:15'
$ws.Range("R131").Value = 'logging iob internet'
$ws.Range("S131").Value = '2024-09-03 20:09:12'
$ws.Range("R132").Value = 'password internet'
$ws.Range("S132").Value = '2024-09-03 20:05:31'
$ws.Range("R133").Value = 'logging iob internet'
$ws.Range("S133").Value = '2024-09-03 20:05:09'
$ws.Range("R134").Value = 'internet'
$ws.Range("S134").Value = '2024-09-03 19:58:18'
$ws.Range("R135").Value = 'login internet invalid'
$ws.Range("S135").Value = '2024-09-03 19:54:49'
$ws.Range("R136").Value = 'login internet invalid'
$ws.Range("S136").Value = '2024-09-03 19:56:17'
$ws.Range("R137").Value = 'corporate internet share'
$ws.Range("S137").Value = '2024-09-03 19:22:58'
$ws.Range("R138").Value = 'login sbi internet personal do not share anyone'
$ws.Range("S138").Value = '2024-09-03 19:17:10'
$ws.Range("R139").Value = 'login internet personal share'
$ws.Range("S139").Value = '2024-09-03 19:13:40'
$ws.Range("R140").Value = 'internet verify it'
$ws.Range("S140").Value = '2024-09-03 19:05:49'
$ws.Range("R141").Value = 'balance your axis'
$ws.Range("S141").Value = '2024-09-03 13:14:06'
$ws.Range("R142").Value = 'lounge'
$ws.Range("S142").Value = '2024-09-03 13:08:08'
$ws.Range("R143").Value = 'balance your axis'
$ws.Range("S143").Value = '2024-09-03 11:21:30'
$ws.Range("R144").Value = 'broker'
$ws.Range("S144").Value = '2024-09-01 22:35:38'

# Move the "Broadband" label from A152 down to the newly added A153 row.
$ws.Range("A152").Value = ""
$ws.Range("A153").Value = "Broadband"

